# Scheduled runner: refresh market-price-derived profit columns across all Sheets.
# Mirrors an upstream price-feed pull: currentAveragePrice/NQ/HQ (H/I/J), the
# LevePrice* columns (K/L), and the resulting LeveProfit* columns (M/N) are
# recomputed per Leve row; a handful of rows also gain/lose their HQ-profit cell
# depending on whether an HQ price is available this cycle.
$wb = $excel.ActiveWorkbook
$updatedCells = 0

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 150
$updatedCells++
$ws.Range("I18").Value = 150
$updatedCells++
$ws.Range("K18").Value = 150
$updatedCells++
$ws.Range("M18").Value = 134
$updatedCells++
$ws.Range("H53").Value = 664.75
$updatedCells++
$ws.Range("I53").Value = 665.6667
$updatedCells++
$ws.Range("J53").Value = 662
$updatedCells++
$ws.Range("K53").Value = 665.6667
$updatedCells++
$ws.Range("L53").Value = 662
$updatedCells++
$ws.Range("M53").Value = -28.66669999999999
$updatedCells++
$ws.Range("N53").Value = -1936
$updatedCells++
$ws.Range("H55").Value = 243.45
$updatedCells++
$ws.Range("I55").Value = 106.75
$updatedCells++
$ws.Range("J55").Value = 334.58334
$updatedCells++
$ws.Range("K55").Value = 106.75
$updatedCells++
$ws.Range("L55").Value = 334.58334
$updatedCells++
$ws.Range("M55").Value = 107.25
$updatedCells++
$ws.Range("N55").Value = -762.58334
$updatedCells++
$ws.Range("H62").Value = 1498.5
$updatedCells++
$ws.Range("I62").Value = 1498.5
$updatedCells++
$ws.Range("K62").Value = 1498.5
$updatedCells++
$ws.Range("M62").Value = -874.5
$updatedCells++
$ws.Range("H65").Value = 1498.5
$updatedCells++
$ws.Range("I65").Value = 1498.5
$updatedCells++
$ws.Range("K65").Value = 7492.5
$updatedCells++
$ws.Range("M65").Value = -4372.5
$updatedCells++
$ws.Range("H80").Value = 279.18182
$updatedCells++
$ws.Range("J80").Value = 290.55554
$updatedCells++
$ws.Range("L80").Value = 871.66662
$updatedCells++
$ws.Range("N80").Value = -2867.66662
$updatedCells++
$ws.Range("H83").Value = 279.18182
$updatedCells++
$ws.Range("J83").Value = 290.55554
$updatedCells++
$ws.Range("L83").Value = 2614.99986
$updatedCells++
$ws.Range("N83").Value = -12598.99986
$updatedCells++
$ws.Range("H103").Value = 1869.375
$updatedCells++
$ws.Range("I103").Value = 0
$updatedCells++
$ws.Range("K103").Value = 0
$updatedCells++
$ws.Range("M103").Value = $null
$updatedCells++
$ws.Range("H111").Value = 2294
$updatedCells++
$ws.Range("I111").Value = 2294
$updatedCells++
$ws.Range("K111").Value = 6882
$updatedCells++
$ws.Range("M111").Value = -3815
$updatedCells++
$ws.Range("H138").Value = 3240.0715
$updatedCells++
$ws.Range("I138").Value = 1738.75
$updatedCells++
$ws.Range("J138").Value = 3840.6
$updatedCells++
$ws.Range("K138").Value = 5216.25
$updatedCells++
$ws.Range("L138").Value = 11521.8
$updatedCells++
$ws.Range("M138").Value = -76.25
$updatedCells++
$ws.Range("N138").Value = -21801.8
$updatedCells++

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1126.5714
$updatedCells++
$ws.Range("I2").Value = 1126.5714
$updatedCells++
$ws.Range("K2").Value = 1126.5714
$updatedCells++
$ws.Range("M2").Value = -1013.5714
$updatedCells++
$ws.Range("H34").Value = 7166579
$updatedCells++
$ws.Range("I34").Value = 16682499
$updatedCells++
$ws.Range("K34").Value = 16682499
$updatedCells++
$ws.Range("M34").Value = -16682228
$updatedCells++
$ws.Range("H40").Value = 0
$updatedCells++
$ws.Range("I40").Value = 0
$updatedCells++
$ws.Range("J40").Value = 0
$updatedCells++
$ws.Range("K40").Value = 0
$updatedCells++
$ws.Range("L40").Value = 0
$updatedCells++
$ws.Range("M40").Value = $null
$updatedCells++
$ws.Range("N40").Value = $null
$updatedCells++
$ws.Range("H42").Value = 31993.889
$updatedCells++
$ws.Range("I42").Value = 28789.4
$updatedCells++
$ws.Range("J42").Value = 35999.5
$updatedCells++
$ws.Range("K42").Value = 28789.4
$updatedCells++
$ws.Range("L42").Value = 35999.5
$updatedCells++
$ws.Range("M42").Value = -28303.4
$updatedCells++
$ws.Range("N42").Value = -36971.5
$updatedCells++
$ws.Range("H61").Value = 2655.7932
$updatedCells++
$ws.Range("I61").Value = 2496.2593
$updatedCells++
$ws.Range("K61").Value = 2496.2593
$updatedCells++
$ws.Range("M61").Value = -2284.2593
$updatedCells++
$ws.Range("H102").Value = 6482.6313
$updatedCells++
$ws.Range("I102").Value = 3021.4
$updatedCells++
$ws.Range("J102").Value = 10328.444
$updatedCells++
$ws.Range("K102").Value = 3021.4
$updatedCells++
$ws.Range("L102").Value = 10328.444
$updatedCells++
$ws.Range("M102").Value = -1399.4
$updatedCells++
$ws.Range("N102").Value = -13572.444
$updatedCells++
$ws.Range("H116").Value = 1126.5714
$updatedCells++
$ws.Range("I116").Value = 1126.5714
$updatedCells++
$ws.Range("K116").Value = 1126.5714
$updatedCells++
$ws.Range("M116").Value = 1167.4286
$updatedCells++
$ws.Range("H136").Value = 2655.7932
$updatedCells++
$ws.Range("I136").Value = 2496.2593
$updatedCells++
$ws.Range("K136").Value = 7488.777900000001
$updatedCells++
$ws.Range("M136").Value = -4938.777900000001
$updatedCells++

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 398.33334
$updatedCells++
$ws.Range("I22").Value = 398.33334
$updatedCells++
$ws.Range("K22").Value = 398.33334
$updatedCells++
$ws.Range("M22").Value = -225.33334
$updatedCells++
$ws.Range("H109").Value = 100342
$updatedCells++
$ws.Range("J109").Value = 100342
$updatedCells++
$ws.Range("L109").Value = 100342
$updatedCells++
$ws.Range("N109").Value = -103116
$updatedCells++
$ws.Range("H134").Value = 3866.6428
$updatedCells++
$ws.Range("J134").Value = 4294
$updatedCells++
$ws.Range("L134").Value = 12882
$updatedCells++
$ws.Range("N134").Value = -17952
$updatedCells++

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 707.5454999999999
$updatedCells++
$ws.Range("I22").Value = 667.4
$updatedCells++
$ws.Range("K22").Value = 667.4
$updatedCells++
$ws.Range("M22").Value = -317.4
$updatedCells++
$ws.Range("H105").Value = 1909.4
$updatedCells++
$ws.Range("I105").Value = 1817
$updatedCells++
$ws.Range("K105").Value = 1817
$updatedCells++
$ws.Range("M105").Value = -70
$updatedCells++
$ws.Range("H107").Value = 649.5
$updatedCells++
$ws.Range("I107").Value = 532.6667
$updatedCells++
$ws.Range("J107").Value = 1000
$updatedCells++
$ws.Range("K107").Value = 532.6667
$updatedCells++
$ws.Range("L107").Value = 1000
$updatedCells++
$ws.Range("M107").Value = 1387.3333
$updatedCells++
$ws.Range("N107").Value = -4840
$updatedCells++

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 189
$updatedCells++
$ws.Range("I97").Value = 188.66667
$updatedCells++
$ws.Range("K97").Value = 566.00001
$updatedCells++
$ws.Range("M97").Value = -70.00000999999997
$updatedCells++

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 20555.5
$updatedCells++
$ws.Range("I35").Value = 0
$updatedCells++
$ws.Range("J35").Value = 20555.5
$updatedCells++
$ws.Range("K35").Value = 0
$updatedCells++
$ws.Range("L35").Value = 20555.5
$updatedCells++
$ws.Range("M35").Value = $null
$updatedCells++
$ws.Range("N35").Value = -21151.5
$updatedCells++
$ws.Range("H44").Value = 23999
$updatedCells++
$ws.Range("J44").Value = 23999
$updatedCells++
$ws.Range("L44").Value = 23999
$updatedCells++
$ws.Range("N44").Value = -25191
$updatedCells++
$ws.Range("H70").Value = 5007978
$updatedCells++
$ws.Range("I70").Value = 10006527
$updatedCells++
$ws.Range("K70").Value = 10006527
$updatedCells++
$ws.Range("M70").Value = -10006257
$updatedCells++
$ws.Range("H73").Value = 5007978
$updatedCells++
$ws.Range("I73").Value = 10006527
$updatedCells++
$ws.Range("K73").Value = 10006527
$updatedCells++
$ws.Range("M73").Value = -10005591
$updatedCells++
$ws.Range("H113").Value = 1862.6666
$updatedCells++
$ws.Range("I113").Value = 1862.6666
$updatedCells++
$ws.Range("K113").Value = 1862.6666
$updatedCells++
$ws.Range("M113").Value = 307.3334
$updatedCells++
$ws.Range("H122").Value = 2499.3333
$updatedCells++
$ws.Range("I122").Value = 2499.3333
$updatedCells++
$ws.Range("K122").Value = 7497.999899999999
$updatedCells++
$ws.Range("M122").Value = -5047.999899999999
$updatedCells++
$ws.Range("H132").Value = 2956
$updatedCells++
$ws.Range("I132").Value = 2615.5
$updatedCells++
$ws.Range("K132").Value = 7846.5
$updatedCells++
$ws.Range("M132").Value = -5316.5
$updatedCells++

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1234.4546
$updatedCells++
$ws.Range("I22").Value = 1050
$updatedCells++
$ws.Range("J22").Value = 1388.1666
$updatedCells++
$ws.Range("K22").Value = 1050
$updatedCells++
$ws.Range("L22").Value = 1388.1666
$updatedCells++
$ws.Range("M22").Value = -755
$updatedCells++
$ws.Range("N22").Value = -1978.1666
$updatedCells++
$ws.Range("H27").Value = 1234.4546
$updatedCells++
$ws.Range("I27").Value = 1050
$updatedCells++
$ws.Range("J27").Value = 1388.1666
$updatedCells++
$ws.Range("K27").Value = 1050
$updatedCells++
$ws.Range("L27").Value = 1388.1666
$updatedCells++
$ws.Range("M27").Value = -943
$updatedCells++
$ws.Range("N27").Value = -1602.1666
$updatedCells++
$ws.Range("H55").Value = 444.625
$updatedCells++
$ws.Range("I55").Value = 324
$updatedCells++
$ws.Range("K55").Value = 324
$updatedCells++
$ws.Range("M55").Value = -151
$updatedCells++
$ws.Range("H68").Value = 2931.25
$updatedCells++
$ws.Range("I68").Value = 2850
$updatedCells++
$ws.Range("K68").Value = 2850
$updatedCells++
$ws.Range("M68").Value = -2101
$updatedCells++
$ws.Range("H71").Value = 2931.25
$updatedCells++
$ws.Range("I71").Value = 2850
$updatedCells++
$ws.Range("K71").Value = 14250
$updatedCells++
$ws.Range("M71").Value = -10506
$updatedCells++

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 627.5
$updatedCells++
$ws.Range("J9").Value = 905
$updatedCells++
$ws.Range("L9").Value = 905
$updatedCells++
$ws.Range("N9").Value = -1185
$updatedCells++
$ws.Range("H40").Value = 37514
$updatedCells++
$ws.Range("I40").Value = 40000
$updatedCells++
$ws.Range("J40").Value = 35028
$updatedCells++
$ws.Range("K40").Value = 40000
$updatedCells++
$ws.Range("L40").Value = 35028
$updatedCells++
$ws.Range("M40").Value = -39851
$updatedCells++
$ws.Range("N40").Value = -35326
$updatedCells++
$ws.Range("H42").Value = 46474.5
$updatedCells++
$ws.Range("I42").Value = 44632.668
$updatedCells++
$ws.Range("K42").Value = 44632.668
$updatedCells++
$ws.Range("M42").Value = -44254.668
$updatedCells++
$ws.Range("H81").Value = 1161.6666
$updatedCells++
$ws.Range("I81").Value = 1094
$updatedCells++
$ws.Range("K81").Value = 2188
$updatedCells++
$ws.Range("M81").Value = -1127
$updatedCells++
$ws.Range("H84").Value = 1161.6666
$updatedCells++
$ws.Range("I84").Value = 1094
$updatedCells++
$ws.Range("K84").Value = 10940
$updatedCells++
$ws.Range("M84").Value = -5636
$updatedCells++
$ws.Range("H100").Value = 441.63635
$updatedCells++
$ws.Range("I100").Value = 448
$updatedCells++
$ws.Range("K100").Value = 896
$updatedCells++
$ws.Range("M100").Value = -355
$updatedCells++
$ws.Range("H108").Value = 90000
$updatedCells++
$ws.Range("J108").Value = 90000
$updatedCells++
$ws.Range("L108").Value = 90000
$updatedCells++
$ws.Range("N108").Value = -97680
$updatedCells++
$ws.Range("H136").Value = 3780.182
$updatedCells++
$ws.Range("I136").Value = 3704.25
$updatedCells++
$ws.Range("J136").Value = 3982.6667
$updatedCells++
$ws.Range("K136").Value = 11112.75
$updatedCells++
$ws.Range("L136").Value = 11948.0001
$updatedCells++
$ws.Range("M136").Value = -8562.75
$updatedCells++
$ws.Range("N136").Value = -17048.0001
$updatedCells++

Write-Host "Updated $updatedCells cells across $($wb.Worksheets.Count) sheets"
